$d = $word.ActiveDocument

# --- Anchor: last paragraph in the body ("BILAGA 1 - Fridlysta arter") ---
$anchor = $d.Paragraphs.Last
$prevRange = $anchor.Range

# =========================================================================
# Phase 1: insert every new paragraph as plain (non-italic) text. We do all
# typing first and only apply character formatting afterwards (Phase 2) --
# this engine treats the most recently assigned Font.Italic as a "sticky"
# typing default for whatever is inserted next (even in a later paragraph),
# so interleaving inserts with formatting would bleed italics forward.
# =========================================================================
# Paragraph 1/13 (style=Heading1)
$prevRange.InsertParagraphAfter()
$para1 = $d.Paragraphs.Last
$para1.Style = "Heading1"
$nr = $para1.Range
$nr.Collapse(0)
$nr.InsertAfter("Knärot – ekologi samt krav på livsmiljön")
$prevRange = $para1.Range

# Paragraph 2/13
$prevRange.InsertParagraphAfter()
$para2 = $d.Paragraphs.Last
$para2.Style = "Normal"
$nr = $para2.Range
$nr.Collapse(0)
$nr.InsertAfter("Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).")
$prevRange = $para2.Range

# Paragraph 3/13
$prevRange.InsertParagraphAfter()
$para3 = $d.Paragraphs.Last
$para3.Style = "Normal"
$nr = $para3.Range
$nr.Collapse(0)
$nr.InsertAfter("Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”")
$prevRange = $para3.Range

# Paragraph 4/13
$prevRange.InsertParagraphAfter()
$para4 = $d.Paragraphs.Last
$para4.Style = "Normal"
$nr = $para4.Range
$nr.Collapse(0)
$nr.InsertAfter("Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”")
$prevRange = $para4.Range

# Paragraph 5/13
$prevRange.InsertParagraphAfter()
$para5 = $d.Paragraphs.Last
$para5.Style = "Normal"
$nr = $para5.Range
$nr.Collapse(0)
$nr.InsertAfter("En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).")
$prevRange = $para5.Range

# Paragraph 6/13
$prevRange.InsertParagraphAfter()
$para6 = $d.Paragraphs.Last
$para6.Style = "Normal"
$nr = $para6.Range
$nr.Collapse(0)
$nr.InsertAfter("Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).")
$prevRange = $para6.Range

# Paragraph 7/13 (style=Heading2)
$prevRange.InsertParagraphAfter()
$para7 = $d.Paragraphs.Last
$para7.Style = "Heading2"
$nr = $para7.Range
$nr.Collapse(0)
$nr.InsertAfter("Referenser - knärot")
$prevRange = $para7.Range

# Paragraph 8/13
$prevRange.InsertParagraphAfter()
$para8 = $d.Paragraphs.Last
$para8.Style = "Normal"
$nr = $para8.Range
$nr.Collapse(0)
$nr.InsertAfter("de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025")
$prevRange = $para8.Range

# Paragraph 9/13
$prevRange.InsertParagraphAfter()
$para9 = $d.Paragraphs.Last
$para9.Style = "Normal"
$nr = $para9.Range
$nr.Collapse(0)
$nr.InsertAfter("Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 ")
$prevRange = $para9.Range

# Paragraph 10/13
$prevRange.InsertParagraphAfter()
$para10 = $d.Paragraphs.Last
$para10.Style = "Normal"
$nr = $para10.Range
$nr.Collapse(0)
$nr.InsertAfter("Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853")
$prevRange = $para10.Range

# Paragraph 11/13
$prevRange.InsertParagraphAfter()
$para11 = $d.Paragraphs.Last
$para11.Style = "Normal"
$nr = $para11.Range
$nr.Collapse(0)
$nr.InsertAfter("Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.")
$prevRange = $para11.Range

# Paragraph 12/13
$prevRange.InsertParagraphAfter()
$para12 = $d.Paragraphs.Last
$para12.Style = "Normal"
$nr = $para12.Range
$nr.Collapse(0)
$nr.InsertAfter("Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/")
$prevRange = $para12.Range

# Paragraph 13/13
$prevRange.InsertParagraphAfter()
$para13 = $d.Paragraphs.Last
$para13.Style = "Normal"
$nr = $para13.Range
$nr.Collapse(0)
$nr.InsertAfter("SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala ")
$prevRange = $para13.Range

# =========================================================================
# Phase 2: re-find each italic run exact text within its own paragraph
# range and flip Font.Italic on just that sub-range. Every new paragraph
# text is already final at this point, so this cannot bleed into anything.
# =========================================================================
# Paragraph 3/13 italics
$fr = $para3.Range.Duplicate
$fr.Find.Execute("“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true
$fr = $para3.Range.Duplicate
$fr.Find.Execute("“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true
$fr = $para3.Range.Duplicate
$fr.Find.Execute("“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true

# Paragraph 4/13 italics
$fr = $para4.Range.Duplicate
$fr.Find.Execute("“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true

# Paragraph 8/13 italics
$fr = $para8.Range.Duplicate
$fr.Find.Execute("Short-term response of the herbaceous layer within leave patches after harvest. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true

# Paragraph 9/13 italics
$fr = $para9.Range.Duplicate
$fr.Find.Execute("Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true

# Paragraph 10/13 italics
$fr = $para10.Range.Duplicate
$fr.Find.Execute("Interactive effects of drought and edge exposure on old-growth forest understory species. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true

# Paragraph 11/13 italics
$fr = $para11.Range.Duplicate
$fr.Find.Execute("Biological legacies buffer local species extinction after logging. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true

# Paragraph 12/13 italics
$fr = $para12.Range.Duplicate
$fr.Find.Execute("Vägledning för hänsyn till knärot. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true

# Paragraph 13/13 italics
$fr = $para13.Range.Duplicate
$fr.Find.Execute("Artfaktablad. Naturvård – artfakta. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Font.Italic = $true

# --- Update the date stamp in the header (2023-09-13 -> 2023-09-15) ---
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $hr = $hdr.Range.Duplicate
            $hr.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}

Write-Output "done"